$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new JSON result strings to D7 and E7 (new "model for result")
$ws.Range("D7").Value = '{"bookingid":1, "userid":1,"doctorid":2, "time":"11:30", "date": "24/12/2023"}'
$ws.Range("E7").Value = '{"status": "Appoinment Successful"}'

# Match the style (alignment) used by the other cells in columns D/E (style index 1:
# centered, vertical-centered, wrap text)
$ws.Range("D7:E7").HorizontalAlignment = -4108
$ws.Range("D7:E7").VerticalAlignment = -4108
$ws.Range("D7:E7").WrapText = $true

# Row 7 grows to accommodate the wrapped text
$ws.Rows.Item(7).RowHeight = 28.8

# Update the view: zoom and selection moved to D7
$ws.Activate()
$excel.ActiveWindow.Zoom = 111
$ws.Range("D7").Select()
